$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vt = [char]11

$cellData = @{
    "1_1" = "53 x 10" + $vt + "  1    0" + $vt + "  ----" + $vt + "5|    |" + $vt + "3|    |"
    "1_2" = "81 x 11" + $vt + "  1    1" + $vt + "  ----" + $vt + "8|    |" + $vt + "1|    |"
    "1_3" = "36 x 20" + $vt + "  2    0" + $vt + "  ----" + $vt + "3|    |" + $vt + "6|    |"
    "2_1" = "68 x 78" + $vt + "  7    8" + $vt + "  ----" + $vt + "6|    |" + $vt + "8|    |"
    "2_2" = "70 x 13" + $vt + "  1    3" + $vt + "  ----" + $vt + "7|    |" + $vt + "0|    |"
    "2_3" = "46 x 55" + $vt + "  5    5" + $vt + "  ----" + $vt + "4|    |" + $vt + "6|    |"
    "3_1" = "61 x 92" + $vt + "  9    2" + $vt + "  ----" + $vt + "6|    |" + $vt + "1|    |"
    "3_2" = "90 x 61" + $vt + "  6    1" + $vt + "  ----" + $vt + "9|    |" + $vt + "0|    |"
    "3_3" = "23 x 20" + $vt + "  2    0" + $vt + "  ----" + $vt + "2|    |" + $vt + "3|    |"
    "4_1" = "67 x 28" + $vt + "  2    8" + $vt + "  ----" + $vt + "6|    |" + $vt + "7|    |"
    "4_2" = "23 x 11" + $vt + "  1    1" + $vt + "  ----" + $vt + "2|    |" + $vt + "3|    |"
    "4_3" = "11 x 87" + $vt + "  8    7" + $vt + "  ----" + $vt + "1|    |" + $vt + "1|    |"
    "5_1" = "61 x 19" + $vt + "  1    9" + $vt + "  ----" + $vt + "6|    |" + $vt + "1|    |"
    "5_2" = "43 x 33" + $vt + "  3    3" + $vt + "  ----" + $vt + "4|    |" + $vt + "3|    |"
    "5_3" = "16 x 61" + $vt + "  6    1" + $vt + "  ----" + $vt + "1|    |" + $vt + "6|    |"
}

for ($r = 1; $r -le 5; $r++) {
    for ($c = 1; $c -le 3; $c++) {
        $cell = $t.Cell($r, $c)
        $key = "$r" + "_" + "$c"
        $cell.Range.Text = $cellData[$key]
    }
}

Write-Host "Done updating lattice multiplication table"